# Auto-generated edit script applying numeric corrections to leve-profit
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4394.727
$ws.Range("J17").Value = 5191.5713
$ws.Range("L17").Value = 15574.7139
$ws.Range("N17").Value = -15910.7139
$ws.Range("H112").Value = 2509.1853
$ws.Range("H113").Value = 4761
$ws.Range("J113").Value = 4891.75
$ws.Range("L113").Value = 4891.75
$ws.Range("N113").Value = -11399.75
$ws.Range("H129").Value = 2923
$ws.Range("I129").Value = 1137.8
$ws.Range("J129").Value = 3734.4546
$ws.Range("K129").Value = 3413.4
$ws.Range("L129").Value = 11203.3638
$ws.Range("M129").Value = 1586.6
$ws.Range("N129").Value = -21203.3638
$ws.Range("H132").Value = 2101.125
$ws.Range("I132").Value = 2083.7273
$ws.Range("K132").Value = 6251.1819
$ws.Range("M132").Value = -3721.1819
$ws.Range("H138").Value = 1763.3462
$ws.Range("I138").Value = 769.4
$ws.Range("K138").Value = 2308.2
$ws.Range("M138").Value = 2831.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 15725
$ws.Range("I14").Value = 19300
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 19300
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -19125
$ws.Range("N14").Value = -5350
$ws.Range("H61").Value = 1941.4286
$ws.Range("I61").Value = 1941.4286
$ws.Range("K61").Value = 1941.4286
$ws.Range("M61").Value = -1729.4286
$ws.Range("H122").Value = 2349.0417
$ws.Range("I122").Value = 2290.318
$ws.Range("K122").Value = 6870.954000000001
$ws.Range("M122").Value = -4420.954000000001
$ws.Range("H136").Value = 1941.4286
$ws.Range("I136").Value = 1941.4286
$ws.Range("K136").Value = 5824.2858
$ws.Range("M136").Value = -3274.2858
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2643.75
$ws.Range("I5").Value = 175
$ws.Range("K5").Value = 175
$ws.Range("M5").Value = -62
$ws.Range("H20").Value = 1346.3077
$ws.Range("I20").Value = 648.8570999999999
$ws.Range("J20").Value = 2160
$ws.Range("K20").Value = 648.8570999999999
$ws.Range("L20").Value = 2160
$ws.Range("M20").Value = -401.8570999999999
$ws.Range("N20").Value = -2654
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620
$ws.Range("H99").Value = 2200
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5396

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 264
$ws.Range("I10").Value = 141.33333
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 141.33333
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -2.333329999999989
$ws.Range("N10").Value = -1278
$ws.Range("H59").Value = 33720.11
$ws.Range("J59").Value = 33720.11
$ws.Range("L59").Value = 33720.11
$ws.Range("N59").Value = -36010.11
$ws.Range("H122").Value = 905.6667
$ws.Range("I122").Value = 698.2857
$ws.Range("K122").Value = 2094.8571
$ws.Range("M122").Value = 355.1428999999998
$ws.Range("H139").Value = 60000
$ws.Range("I139").Value = 60000
$ws.Range("K139").Value = 60000
$ws.Range("M139").Value = -54860
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null
$ws.Range("H141").Value = 87373.75
$ws.Range("J141").Value = 96998.57000000001
$ws.Range("L141").Value = 96998.57000000001
$ws.Range("N141").Value = -107358.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1375.5
$ws.Range("I11").Value = 1567.3334
$ws.Range("J11").Value = 800
$ws.Range("K11").Value = 4702.0002
$ws.Range("L11").Value = 2400
$ws.Range("M11").Value = -4562.0002
$ws.Range("N11").Value = -2680
$ws.Range("H34").Value = 1862
$ws.Range("I34").Value = 816
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 2448
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = -2364
$ws.Range("N34").Value = -15168
$ws.Range("H39").Value = 831.6667
$ws.Range("J39").Value = 997.5
$ws.Range("L39").Value = 2992.5
$ws.Range("N39").Value = -3580.5
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = $null
$ws.Range("N92").Value = $null
$ws.Range("H119").Value = 129
$ws.Range("I119").Value = 129
$ws.Range("K119").Value = 387
$ws.Range("M119").Value = 4451
$ws.Range("H120").Value = 13632.5
$ws.Range("J120").Value = 23332.666
$ws.Range("L120").Value = 69997.99800000001
$ws.Range("N120").Value = -79673.99800000001
$ws.Range("H121").Value = 657
$ws.Range("J121").Value = 966
$ws.Range("L121").Value = 2898
$ws.Range("N121").Value = -5518
$ws.Range("H122").Value = 1490
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = $null
$ws.Range("N123").Value = $null
$ws.Range("H137").Value = 3193
$ws.Range("I137").Value = 2664.5
$ws.Range("J137").Value = 4250
$ws.Range("K137").Value = 7993.5
$ws.Range("L137").Value = 12750
$ws.Range("M137").Value = -2893.5
$ws.Range("N137").Value = -22950

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 539
$ws.Range("I2").Value = 532.3333
$ws.Range("J2").Value = 549
$ws.Range("K2").Value = 532.3333
$ws.Range("L2").Value = 549
$ws.Range("M2").Value = -419.3333
$ws.Range("N2").Value = -775
$ws.Range("H126").Value = 4878.5
$ws.Range("I126").Value = 3166.6667
$ws.Range("J126").Value = 10014
$ws.Range("K126").Value = 9500.000100000001
$ws.Range("L126").Value = 30042
$ws.Range("M126").Value = -7030.000100000001
$ws.Range("N126").Value = -34982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9937.714
$ws.Range("I40").Value = 9913
$ws.Range("J40").Value = 9999.5
$ws.Range("K40").Value = 9913
$ws.Range("L40").Value = 9999.5
$ws.Range("M40").Value = -9777
$ws.Range("N40").Value = -10271.5
$ws.Range("H46").Value = 2781.0454
$ws.Range("I46").Value = 2511.1177
$ws.Range("K46").Value = 2511.1177
$ws.Range("M46").Value = -2323.1177
$ws.Range("H123").Value = 78332.664
$ws.Range("J123").Value = 78332.664
$ws.Range("L123").Value = 78332.664
$ws.Range("N123").Value = -88132.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = $null
$ws.Range("H44").Value = 25010
$ws.Range("I44").Value = 25010
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 25010
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -24456
$ws.Range("N44").Value = $null
$ws.Range("H69").Value = 8542.25
$ws.Range("J69").Value = 8542.25
$ws.Range("L69").Value = 8542.25
$ws.Range("N69").Value = -10040.25
$ws.Range("H72").Value = 8542.25
$ws.Range("J72").Value = 8542.25
$ws.Range("L72").Value = 25626.75
$ws.Range("N72").Value = -33114.75
$ws.Range("H81").Value = 218.33333
$ws.Range("J81").Value = 175
$ws.Range("L81").Value = 350
$ws.Range("N81").Value = -2472
$ws.Range("H84").Value = 218.33333
$ws.Range("J84").Value = 175
$ws.Range("L84").Value = 1750
$ws.Range("N84").Value = -12358
$ws.Range("H122").Value = 1904.4445
$ws.Range("I122").Value = 1904.4445
$ws.Range("K122").Value = 5713.333500000001
$ws.Range("M122").Value = -3263.333500000001
$ws.Range("H136").Value = 1858.2106
$ws.Range("I136").Value = 1980.8
$ws.Range("K136").Value = 5942.4
$ws.Range("M136").Value = -3392.4
$ws.Range("H138").Value = 100000
$ws.Range("I138").Value = 100000
$ws.Range("K138").Value = 100000
$ws.Range("M138").Value = -94860
